$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: columns A-J ("<Name>_old") become "<Name>_FV2410",
#    column K ("diff") is left untouched, columns L-U ("<Name>_new") become
#    "<Name>_FV2504".
$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$headersFV2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headersFV2410.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2410[$i]
}

for ($i = 0; $i -lt $headersFV2504.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2504[$i]
}

# 2) Freeze the header row (split below row 1, pane anchored at A2).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the used range into an Excel Table ("Table1") with autofilter.
$rng = $ws.Range("A1:U71")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
